# Update "想去人数" (interest count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 646
$ws1.Range("F4").Value  = 210
$ws1.Range("F6").Value  = 9796
$ws1.Range("F7").Value  = 884
$ws1.Range("F9").Value  = 1236
$ws1.Range("F10").Value = 3383
$ws1.Range("F11").Value = 169
$ws1.Range("F12").Value = 113
$ws1.Range("F13").Value = 38
$ws1.Range("F16").Value = 535
$ws1.Range("F18").Value = 265
$ws1.Range("F19").Value = 1430

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 18

# --- Sheet "全部类型" (all types, aggregate) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 18
$ws4.Range("F4").Value  = 646
$ws4.Range("F5").Value  = 210
$ws4.Range("F7").Value  = 9796
$ws4.Range("F8").Value  = 884
$ws4.Range("F10").Value = 1236
$ws4.Range("F11").Value = 3383
$ws4.Range("F12").Value = 169
$ws4.Range("F13").Value = 113
$ws4.Range("F14").Value = 38
$ws4.Range("F17").Value = 535
$ws4.Range("F19").Value = 265
$ws4.Range("F20").Value = 1430
